$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Form Pemesanan: mark tables (meja) 1-6 as booked ("dipesan") in column B
$ws.Range("B1:B6").Value = 1
